$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: updated timestamp
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty -> "Alvearie Team"
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Old row 11 was a duplicate "Contact" / "No display for ContactDetail" row - delete it,
# shifting rows 12-15 up to 11-14
$ws.Rows("11").Delete()
